$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1698
$ws.Range("E2").Value = 112
$ws.Range("F2").Value = 112
$ws.Range("G2").Value = -29
$ws.Range("H2").Value = -26
$ws.Range("I2").Value = -26
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1630
$ws.Range("L2").Value = 710
$ws.Range("M2").Value = 920
$ws.Range("N2").Value = 920
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 140
$ws.Range("Q2").Value = 194
$ws.Range("R2").Value = -104
$ws.Range("S2").Value = -18
$ws.Range("T2").Value = 30
$ws.Range("U2").Value = 164
$ws.Range("V2").Value = 387
$ws.Range("W2").Value = 6.58
$ws.Range("X2").Value = -1.51
$ws.Range("Y2").Value = -2.74
$ws.Range("Z2").Value = -1.63
$ws.Range("AA2").Value = 77.12
$ws.Range("AB2").Value = 576.37
$ws.Range("AC2").Value = -106
$ws.Range("AD2").Value = -42.68
$ws.Range("AE2").Value = 3802
$ws.Range("AF2").Value = 1.19
$ws.Range("AG2").Value = 60
$ws.Range("AH2").Value = 1.32
$ws.Range("AI2").Value = -56.53
$ws.Range("AJ2").Value = 24200000

# Row 3
$ws.Range("D3").Value = 1721
$ws.Range("E3").Value = 178
$ws.Range("F3").Value = 178
$ws.Range("G3").Value = 225
$ws.Range("H3").Value = 154
$ws.Range("I3").Value = 154
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1706
$ws.Range("L3").Value = 670
$ws.Range("M3").Value = 1036
$ws.Range("N3").Value = 1036
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 140
$ws.Range("Q3").Value = 42
$ws.Range("R3").Value = -21
$ws.Range("S3").Value = -7
$ws.Range("T3").Value = 22
$ws.Range("U3").Value = 20
$ws.Range("V3").Value = 405
$ws.Range("W3").Value = 10.35
$ws.Range("X3").Value = 8.970000000000001
$ws.Range("Y3").Value = 15.78
$ws.Range("Z3").Value = 9.25
$ws.Range("AA3").Value = 64.68000000000001
$ws.Range("AB3").Value = 668.63
$ws.Range("AC3").Value = 638
$ws.Range("AD3").Value = 7.41
$ws.Range("AE3").Value = 4327
$ws.Range("AF3").Value = 1.09
$ws.Range("AG3").Value = 60
$ws.Range("AH3").Value = 1.27
$ws.Range("AI3").Value = 9.32
$ws.Range("AJ3").Value = 24200000

# Row 4
$ws.Range("D4").Value = 2008
$ws.Range("E4").Value = 276
$ws.Range("F4").Value = 276
$ws.Range("G4").Value = 293
$ws.Range("H4").Value = 216
$ws.Range("I4").Value = 219
$ws.Range("J4").Value = -2
$ws.Range("K4").Value = 1868
$ws.Range("L4").Value = 637
$ws.Range("M4").Value = 1231
$ws.Range("N4").Value = 1221
$ws.Range("O4").Value = 9
$ws.Range("P4").Value = 140
$ws.Range("Q4").Value = 220
$ws.Range("R4").Value = -70
$ws.Range("S4").Value = -138
$ws.Range("T4").Value = 74
$ws.Range("U4").Value = 146
$ws.Range("V4").Value = 287
$ws.Range("W4").Value = 13.76
$ws.Range("X4").Value = 10.76
$ws.Range("Y4").Value = 19.37
$ws.Range("Z4").Value = 12.09
$ws.Range("AA4").Value = 51.78
$ws.Range("AB4").Value = 812.55
$ws.Range("AC4").Value = 903
$ws.Range("AD4").Value = 10.49
$ws.Range("AE4").Value = 5169
$ws.Range("AF4").Value = 1.83
$ws.Range("AG4").Value = 90
$ws.Range("AH4").Value = 0.95
$ws.Range("AI4").Value = 9.73
$ws.Range("AJ4").Value = 24200000

# Row 5
$ws.Range("D5").Value = 2114
$ws.Range("E5").Value = 237
$ws.Range("F5").Value = 237
$ws.Range("G5").Value = 146
$ws.Range("H5").Value = 87
$ws.Range("I5").Value = 90
$ws.Range("J5").Value = -3
$ws.Range("K5").Value = 1918
$ws.Range("L5").Value = 668
$ws.Range("M5").Value = 1250
$ws.Range("N5").Value = 1250
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 140
$ws.Range("Q5").Value = 212
$ws.Range("R5").Value = -51
$ws.Range("S5").Value = -23
$ws.Range("T5").Value = 32
$ws.Range("U5").Value = 180
$ws.Range("V5").Value = 306
$ws.Range("W5").Value = 11.2
$ws.Range("X5").Value = 4.13
$ws.Range("Y5").Value = 7.27
$ws.Range("Z5").Value = 4.61
$ws.Range("AA5").Value = 53.43
$ws.Range("AB5").Value = 858.45
$ws.Range("AC5").Value = 371
$ws.Range("AD5").Value = 18.69
$ws.Range("AE5").Value = 5286
$ws.Range("AF5").Value = 1.31
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 1.44
$ws.Range("AI5").Value = 26.32
$ws.Range("AJ5").Value = 24200000

# Row 6
$ws.Range("D6").Value = 3231
$ws.Range("E6").Value = 215
$ws.Range("F6").Value = 215
$ws.Range("G6").Value = 219
$ws.Range("H6").Value = 149
$ws.Range("I6").Value = 153
$ws.Range("K6").Value = 2947
$ws.Range("L6").Value = 1566
$ws.Range("M6").Value = 1381
$ws.Range("N6").Value = 1369
$ws.Range("P6").Value = 140
$ws.Range("Q6").Value = -80
$ws.Range("R6").Value = -366
$ws.Range("S6").Value = 337
$ws.Range("T6").Value = 49
$ws.Range("U6").Value = -129
$ws.Range("V6").Value = 985
$ws.Range("W6").Value = 6.66
$ws.Range("X6").Value = 4.63
$ws.Range("Y6").Value = 11.66
$ws.Range("Z6").Value = 6.14
$ws.Range("AA6").Value = 113.44
$ws.Range("AB6").Value = 942.51
$ws.Range("AC6").Value = 631
$ws.Range("AD6").Value = 10.04
$ws.Range("AE6").Value = 5789
$ws.Range("AF6").Value = 1.1
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 1.58
$ws.Range("AI6").Value = 15.49
$ws.Range("AJ6").Value = 24200000

# Row 7
$ws.Range("D7").Value = 4658
$ws.Range("E7").Value = 407
$ws.Range("G7").Value = 372
$ws.Range("H7").Value = 264
$ws.Range("I7").Value = 270
$ws.Range("K7").Value = 3402
$ws.Range("L7").Value = 1780
$ws.Range("M7").Value = 1621
$ws.Range("N7").Value = 1615
$ws.Range("P7").Value = 140
$ws.Range("Q7").Value = -105
$ws.Range("R7").Value = 205
$ws.Range("S7").Value = -21
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = -105
$ws.Range("W7").Value = 8.74
$ws.Range("X7").Value = 5.67
$ws.Range("Y7").Value = 18.09
$ws.Range("Z7").Value = 8.32
$ws.Range("AA7").Value = 109.81
$ws.Range("AC7").Value = 1116
$ws.Range("AD7").Value = 6.59
$ws.Range("AE7").Value = 6827
$ws.Range("AF7").Value = 1.08
$ws.Range("AG7").Value = 100
$ws.Range("AH7").Value = 1.36
$ws.Range("AI7").Value = 8.960000000000001

# Row 8
$ws.Range("D8").Value = 5245
$ws.Range("E8").Value = 506
$ws.Range("G8").Value = 472
$ws.Range("H8").Value = 377
$ws.Range("I8").Value = 386
$ws.Range("K8").Value = 3798
$ws.Range("L8").Value = 1823
$ws.Range("M8").Value = 1975
$ws.Range("N8").Value = 1978
$ws.Range("P8").Value = 140
$ws.Range("Q8").Value = 233
$ws.Range("R8").Value = 8
$ws.Range("S8").Value = -68
$ws.Range("T8").Value = 0
$ws.Range("U8").Value = 233
$ws.Range("W8").Value = 9.65
$ws.Range("X8").Value = 7.19
$ws.Range("Y8").Value = 21.49
$ws.Range("Z8").Value = 10.47
$ws.Range("AA8").Value = 92.3
$ws.Range("AC8").Value = 1595
$ws.Range("AD8").Value = 4.61
$ws.Range("AE8").Value = 8362
$ws.Range("AF8").Value = 0.88
$ws.Range("AG8").Value = 100
$ws.Range("AH8").Value = 1.36
$ws.Range("AI8").Value = 6.27

# Row 9
$ws.Range("D9").Value = 5763
$ws.Range("E9").Value = 581
$ws.Range("G9").Value = 549
$ws.Range("H9").Value = 439
$ws.Range("I9").Value = 449
$ws.Range("K9").Value = 4291
$ws.Range("L9").Value = 1900
$ws.Range("M9").Value = 2390
$ws.Range("N9").Value = 2402
$ws.Range("P9").Value = 140
$ws.Range("Q9").Value = 310
$ws.Range("R9").Value = 8
$ws.Range("S9").Value = -23
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 310
$ws.Range("W9").Value = 10.08
$ws.Range("X9").Value = 7.62
$ws.Range("Y9").Value = 20.5
$ws.Range("Z9").Value = 10.85
$ws.Range("AA9").Value = 79.5
$ws.Range("AC9").Value = 1855
$ws.Range("AD9").Value = 3.96
$ws.Range("AE9").Value = 10154
$ws.Range("AF9").Value = 0.72
$ws.Range("AG9").Value = 100
$ws.Range("AH9").Value = 1.36
$ws.Range("AI9").Value = 5.39
